# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 13:20"

# Madrid (row 4) - new totals
$ws.Range("B4").Value = 32155
$ws.Range("C4").Value = 12400
$ws.Range("D4").Value = 15580
$ws.Range("E4").Value = 4175

# Navarra (row 9) - new totals
$ws.Range("B9").Value = 2682
$ws.Range("C9").Value = 278
$ws.Range("D9").Value = 2263
$ws.Range("E9").Value = 141

# La Rioja (row 13) - new totals
$ws.Range("B13").Value = 2083
$ws.Range("C13").Value = 641
$ws.Range("D13").Value = 1341
$ws.Range("E13").Value = 101

# Cantabria overtakes Sevilla in ranking (sheet sorted descending by Casos totales),
# so the two rows swap places. Row 23 becomes Cantabria (with its updated, higher
# total), row 24 becomes Sevilla (keeping its previous, now-lower-ranked totals).
$ws.Range("A23").Value = "Cantabria"
$ws.Range("B23").Value = 1268
$ws.Range("C23").Value = 60
$ws.Range("D23").Value = 1148
$ws.Range("E23").Value = 60

$ws.Range("A24").Value = "Sevilla"
$ws.Range("B24").Value = 1215
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 1140
$ws.Range("E24").Value = 55

# Murcia (row 29) - small correction to Recuperados
$ws.Range("D29").Value = 997
